$wb = $excel.ActiveWorkbook

# --- Sheet "Inventario": update Cantidad values (column B) ---
$inv = $wb.Worksheets.Item("Inventario")
$inv.Range("B2").Value = 45
$inv.Range("B3").Value = 26
$inv.Range("B4").Value = 500
$inv.Range("B5").Value = 49
$inv.Range("B6").Value = 20
$inv.Range("B7").Value = 7

# --- Sheet "Ventas": append sales summary rows for the day 2025-10-22 (serial 45952) ---
$ventas = $wb.Worksheets.Item("Ventas")

$newRows = @(
    @{ Fecha = 45952; Producto = "Bolsa de Regalo"; Cantidad = 2;  Unidad = "unidades"; Precio = 3000;  Subtotal = 6000 },
    @{ Fecha = 45952; Producto = "Loción";          Cantidad = 30; Unidad = "gramos";   Precio = 550;   Subtotal = 16500 },
    @{ Fecha = 45952; Producto = "Splash";          Cantidad = 1;  Unidad = "unidades"; Precio = 16000; Subtotal = 16000 },
    @{ Fecha = 45952; Producto = "Chocolatina";      Cantidad = 2;  Unidad = "unidades"; Precio = 6000;  Subtotal = 12000 },
    @{ Fecha = 45952; Producto = "Loción";          Cantidad = 50; Unidad = "gramos";   Precio = 550;   Subtotal = 27500 },
    @{ Fecha = 45952; Producto = "Chocolatina";      Cantidad = 1;  Unidad = "unidades"; Precio = 6000;  Subtotal = 6000 },
    @{ Fecha = 45952; Producto = "Crema";            Cantidad = 3;  Unidad = "unidades"; Precio = 14500; Subtotal = 43500 },
    @{ Fecha = 45952; Producto = "Camisa";           Cantidad = 1;  Unidad = "unidades"; Precio = 32000; Subtotal = 32000 },
    @{ Fecha = 45952; Producto = "Chocolatina";      Cantidad = 3;  Unidad = "unidades"; Precio = 6000;  Subtotal = 18000 },
    @{ Fecha = 45952; Producto = "Camisa";           Cantidad = 2;  Unidad = "unidades"; Precio = 32000; Subtotal = 64000 }
)

$startRow = 16
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $cellA = $ventas.Cells.Item($r, 1)
    $cellA.Value = $row.Fecha
    $cellA.NumberFormat = "YYYY-MM-DD"

    $ventas.Cells.Item($r, 2).Value = $row.Producto
    $ventas.Cells.Item($r, 3).Value = $row.Cantidad
    $ventas.Cells.Item($r, 4).Value = $row.Unidad
    $ventas.Cells.Item($r, 5).Value = $row.Precio
    $ventas.Cells.Item($r, 6).Value = $row.Subtotal
}
